$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Logistic_Regression)
$ws.Range("B2").Value = 0.85
$ws.Range("C2").Value = 0.82
$ws.Range("D2").Value = 12
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 0.67
$ws.Range("I2").Value = 0.71
$ws.Range("J2").Value = 0.6899999999999999
$ws.Range("L2").Value = 0.79

# Row 3 (Random_Forest)
$ws.Range("B3").Value = 0.99
$ws.Range("C3").Value = 0.8100000000000001
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 42
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 9
$ws.Range("H3").Value = 0.73
$ws.Range("I3").Value = 0.47
$ws.Range("J3").Value = 0.57
$ws.Range("K3").Value = 0.93
$ws.Range("L3").Value = 0.7

# Row 4 (Kernel_SVM)
$ws.Range("B4").Value = 0.88
$ws.Range("C4").Value = 0.82
$ws.Range("D4").Value = 12
$ws.Range("E4").Value = 39
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 0.67
$ws.Range("I4").Value = 0.71
$ws.Range("J4").Value = 0.6899999999999999
$ws.Range("K4").Value = 0.87
$ws.Range("L4").Value = 0.79

# Row 5 (CatBoost)
$ws.Range("B5").Value = 0.95
$ws.Range("D5").Value = 12
$ws.Range("E5").Value = 40
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 0.71
$ws.Range("I5").Value = 0.71
$ws.Range("J5").Value = 0.71
$ws.Range("K5").Value = 0.89
$ws.Range("L5").Value = 0.8

# Row 6 (DNN)
$ws.Range("B6").Value = 0.95
$ws.Range("C6").Value = 0.84
$ws.Range("E6").Value = 11
$ws.Range("F6").Value = 34
$ws.Range("H6").Value = 0.33
$ws.Range("J6").Value = 0.5
$ws.Range("K6").Value = 0.24
$ws.Range("L6").Value = 0.62
